$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, pushing existing rows 79..107 down to 80..108
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with data (columns A,B,C,E,F,G,H,I,R are identical to the
# surrounding rows for this market/category subset, so copy them from row 80,
# then set the row-specific values per the diff).
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = 44636
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = 100112032
$ws.Range("G79").Value = "Zapallo italiano"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 220
$ws.Range("K79").Value = 11000
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = 11455
$ws.Range("N79").Value = "$/caja 60 unidades"
$ws.Range("O79").Value = "Provincia de Huasco"
$ws.Range("P79").Value = 191
$ws.Range("Q79").Value = 60
$ws.Range("R79").Value = "Hortaliza"

# Make sure the date cell uses the same date number format style as the rest of column D
$ws.Range("D79").NumberFormat = $ws.Range("D80").NumberFormat
